$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value would otherwise be auto-parsed as a number by Excel
# (losing the literal text representation, e.g. trailing zeros / sci notation).
# Force them to Text format first so COM stores the exact original string.
foreach ($addr in @('D4', 'D5', 'D6', 'D7', 'D9', 'D10', 'D12', 'D13', 'D14', 'D15', 'D17', 'D19', 'D20', 'D21', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '24.922.89'
$ws.Range("E2").Value = '  -3.79%  '

# Row 3
$ws.Range("D3").Value = '1.637.22'
$ws.Range("E3").Value = '  -5.97%  '

# Row 4
$ws.Range("D4").Value = '0.9974'

# Row 5
$ws.Range("D5").Value = '235.75'
$ws.Range("E5").Value = '  -4.55%  '

# Row 6
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.10%  '

# Row 7
$ws.Range("D7").Value = '0.4717'
$ws.Range("E7").Value = '  -6.68%  '

# Row 8
$ws.Range("E8").Value = '  -5.89%  '

# Row 9
$ws.Range("D9").Value = '0.06012'
$ws.Range("E9").Value = '  -2.68%  '

# Row 10
$ws.Range("D10").Value = '0.07145'
$ws.Range("E10").Value = '  -1.18%  '

# Row 11
$ws.Range("D11").Value = '1.632.24'
$ws.Range("E11").Value = '  -6.31%  '

# Row 12
$ws.Range("D12").Value = '14.77'
$ws.Range("E12").Value = '  -1.95%  '

# Row 13
$ws.Range("D13").Value = '0.6138'
$ws.Range("E13").Value = '  -4.97%  '

# Row 14
$ws.Range("D14").Value = '4.414'
$ws.Range("E14").Value = '  -4.50%  '

# Row 15
$ws.Range("D15").Value = '72.55'
$ws.Range("E15").Value = '  -6.30%  '

# Row 16
$ws.Range("E16").Value = '  -0.06%  '

# Row 17
$ws.Range("D17").Value = '0.9976'
$ws.Range("E17").Value = '  -0.32%  '

# Row 18
$ws.Range("D18").Value = '24.906.41'
$ws.Range("E18").Value = '  -3.91%  '

# Row 19
$ws.Range("D19").Value = '0.000006566'
$ws.Range("E19").Value = '  -3.36%  '

# Row 20
$ws.Range("D20").Value = '11.21'
$ws.Range("E20").Value = '  -5.05%  '

# Row 21
$ws.Range("D21").Value = '4.405'
$ws.Range("E21").Value = '  +3.09%  '

# Row 22
$ws.Range("D22").Value = '1.846.81'
$ws.Range("E22").Value = '  -6.23%  '

# Row 23
$ws.Range("D23").Value = '8.553'
$ws.Range("E23").Value = '  -0.79%  '

# Row 24
$ws.Range("D24").Value = '5.250'
$ws.Range("E24").Value = '  -2.40%  '

# Row 25
$ws.Range("D25").Value = '132.49'
$ws.Range("E25").Value = '  -2.63%  '

# Row 26
$ws.Range("D26").Value = '14.75'
$ws.Range("E26").Value = '  -3.13%  '

# Row 27
$ws.Range("D27").Value = '1.373'
$ws.Range("E27").Value = '  -8.55%  '

# Row 28
$ws.Range("E28").Value = '  -3.03%  '

# Row 29
$ws.Range("D29").Value = '1.649'
$ws.Range("E29").Value = '  -6.50%  '

# Row 30
$ws.Range("D30").Value = '3.722'
$ws.Range("E30").Value = '  -4.74%  '

# Row 31
$ws.Range("D31").Value = '0.07738'
$ws.Range("E31").Value = '  -5.86%  '

# Row 32
$ws.Range("D32").Value = '3.536'
$ws.Range("E32").Value = '  -2.44%  '

# Row 33
$ws.Range("D33").Value = '0.04370'
$ws.Range("E33").Value = '  -6.40%  '

# Row 34
$ws.Range("D34").Value = '0.9993'
$ws.Range("E34").Value = '  -0.11%  '

# Row 35
$ws.Range("D35").Value = '2.596'
$ws.Range("E35").Value = '  -2.27%  '

# Row 36
$ws.Range("D36").Value = '0.9156'
$ws.Range("E36").Value = '  -7.70%  '

# Row 37
$ws.Range("D37").Value = '0.5802'
$ws.Range("E37").Value = '  -6.36%  '

# Row 38
$ws.Range("D38").Value = '2.536'
$ws.Range("E38").Value = '  -7.07%  '

# Row 39
$ws.Range("D39").Value = '0.01551'
$ws.Range("E39").Value = '  -2.96%  '

# Row 40
$ws.Range("D40").Value = '0.9979'
$ws.Range("E40").Value = '  -0.30%  '

# Row 41
$ws.Range("D41").Value = '0.8194'
$ws.Range("E41").Value = '  +8.39%  '

# Row 42
$ws.Range("D42").Value = '1.791'
$ws.Range("E42").Value = '  -6.28%  '

# Row 43
$ws.Range("D43").Value = '97.38'
$ws.Range("E43").Value = '  -1.54%  '

# Row 44
$ws.Range("D44").Value = '0.3697'
$ws.Range("E44").Value = '  -4.14%  '

# Row 45
$ws.Range("D45").Value = '4.732'
$ws.Range("E45").Value = '  -4.98%  '

# Row 46
$ws.Range("D46").Value = '0.1127'
$ws.Range("E46").Value = '  -0.66%  '

# Row 47
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05203'
$ws.Range("E47").Value = '  -0.62%  '

# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '6.071'
$ws.Range("E48").Value = '  -2.92%  '

# Row 49
$ws.Range("D49").Value = '29.44'
$ws.Range("E49").Value = '  -3.74%  '

# Row 50
$ws.Range("D50").Value = '0.9995'
$ws.Range("E50").Value = '  -0.40%  '

# Row 51
$ws.Range("D51").Value = '0.9999'
$ws.Range("E51").Value = '  -0.52%  '
